$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New de-duplicated data, prioritizing rows that already have a
# first/last name filled in so they survive the dedupe pass.
$data = @(
    @("Amy",   "Feldkamp", "afeldkamp@troy.k12.mi.us",      "",              ""),
    @("April", "Kosin",    "adkosin@clarkston.k12.mi.us",   "",              ""),
    @("Elsa",  "Garcia",   "273971@dadeschools.net",        "",              ""),
    @("",      "",         "181884@dadeschools.net",        "",              "teacher"),
    @("",      "",         "222963@dadeschools.net",        "",              ""),
    @("",      "",         "251936@dadeschools.net",        "",              ""),
    @("",      "",         "288058@dadeschools.net",        "this@that.com", "teacher"),
    @("",      "",         "296806@dadeschools.net",        "eat@gmail.com", "teacher"),
    @("",      "",         "204093@dadeschools.net",        "",              "teacher"),
    @("",      "",         "addison@dadeschools.net",       "",              "administrator"),
    @("",      "",         "adsturm@dsdmail.net",           "",              ""),
    @("",      "",         "adtaylor@dsdmail.net",          "",              "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
